$d = $word.ActiveDocument

$d.Content.Find.Execute("Projektideen", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Projektideen`r`nNicht viel", 2)

$d.Save()
